# FixedCrops and sum at founds
#
# Adds a new "Fixo" column (K) to the Cultivo sheet, holding a fixed
# value per crop row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (K1) -----------------------------------------------------
# Clone the formatting of the neighbouring header cell so the new
# header cell matches the rest of row 1, then set its text.
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("K1").Value = "Fixo"

# --- Data (K2:K9) ------------------------------------------------------
# Clone the formatting of a neighbouring data cell across the new
# column's data range before filling in the values.
$ws.Range("H2").Copy() | Out-Null
$ws.Range("K2:K9").PasteSpecial(-4122) | Out-Null

$ws.Range("K2").Value = 10.0
$ws.Range("K3").Value = 0.0
$ws.Range("K4").Value = 0.0
$ws.Range("K5").Value = 0.0
$ws.Range("K6").Value = 0.0
$ws.Range("K7").Value = 0.0
$ws.Range("K8").Value = 0.0
$ws.Range("K9").Value = 0.0
